# The "Förändrad" (Changed) date column (C) for every data row (2-151) is
# updated from 2023-09-13 (serial 45182) to 2023-09-15 (serial 45184).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C151").Value = 45184
